$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ALC = $wb.Worksheets.Item("ALC")
$ALC.Range("H40").Value = 2310.8572
$ALC.Range("I40").Value = 3513.1667
$ALC.Range("J40").Value = 1829.9333
$ALC.Range("K40").Value = 3513.1667
$ALC.Range("L40").Value = 1829.9333
$ALC.Range("M40").Value = -3338.1667
$ALC.Range("N40").Value = -2179.9333
$ALC.Range("H62").Value = 2659.8333
$ALC.Range("I62").Value = 2239.75
$ALC.Range("K62").Value = 2239.75
$ALC.Range("M62").Value = -1615.75
$ALC.Range("H64").Value = 2719.4736
$ALC.Range("J64").Value = 2774.6155
$ALC.Range("L64").Value = 2774.6155
$ALC.Range("N64").Value = -3270.6155
$ALC.Range("H65").Value = 2659.8333
$ALC.Range("I65").Value = 2239.75
$ALC.Range("K65").Value = 11198.75
$ALC.Range("M65").Value = -8078.75
$ALC.Range("H67").Value = 2719.4736
$ALC.Range("J67").Value = 2774.6155
$ALC.Range("L67").Value = 2774.6155
$ALC.Range("N67").Value = -4490.6155
$ALC.Range("H74").Value = 4789.0557
$ALC.Range("I74").Value = 5611.4443
$ALC.Range("J74").Value = 3966.6667
$ALC.Range("K74").Value = 5611.4443
$ALC.Range("L74").Value = 3966.6667
$ALC.Range("M74").Value = -4675.4443
$ALC.Range("N74").Value = -5838.6667
$ALC.Range("H76").Value = 55452.79
$ALC.Range("I76").Value = 65118.938
$ALC.Range("J76").Value = 3900
$ALC.Range("K76").Value = 65118.938
$ALC.Range("L76").Value = 3900
$ALC.Range("M76").Value = -64803.938
$ALC.Range("N76").Value = -4530
$ALC.Range("H77").Value = 4789.0557
$ALC.Range("I77").Value = 5611.4443
$ALC.Range("J77").Value = 3966.6667
$ALC.Range("K77").Value = 28057.2215
$ALC.Range("L77").Value = 19833.3335
$ALC.Range("M77").Value = -23377.2215
$ALC.Range("N77").Value = -29193.3335
$ALC.Range("H79").Value = 55452.79
$ALC.Range("I79").Value = 65118.938
$ALC.Range("J79").Value = 3900
$ALC.Range("K79").Value = 65118.938
$ALC.Range("L79").Value = 3900
$ALC.Range("M79").Value = -64026.938
$ALC.Range("N79").Value = -6084
$ALC.Range("H106").Value = 29472974
$ALC.Range("I106").Value = 73611.42999999999
$ALC.Range("K106").Value = 73611.42999999999
$ALC.Range("M106").Value = -72980.42999999999
$ALC.Range("H138").Value = 2487.433
$ALC.Range("I138").Value = 1350.7059
$ALC.Range("J138").Value = 3100.9048
$ALC.Range("K138").Value = 4052.1177
$ALC.Range("L138").Value = 9302.714399999999
$ALC.Range("M138").Value = 1087.8823
$ALC.Range("N138").Value = -19582.7144

# ---- ARM ----
$ARM = $wb.Worksheets.Item("ARM")
$ARM.Range("H28").Value = 11721.5
$ARM.Range("I28").Value = 754.4
$ARM.Range("J28").Value = 30000
$ARM.Range("K28").Value = 754.4
$ARM.Range("L28").Value = 30000
$ARM.Range("M28").Value = -562.4
$ARM.Range("N28").Value = -30384
$ARM.Range("H45").Value = 17544804
$ARM.Range("I45").Value = 37037804
$ARM.Range("J45").Value = 1102.8
$ARM.Range("K45").Value = 37037804
$ARM.Range("L45").Value = 1102.8
$ARM.Range("M45").Value = -37037427
$ARM.Range("N45").Value = -1856.8
$ARM.Range("H61").Value = 1702.5526
$ARM.Range("I61").Value = 1413.5927
$ARM.Range("K61").Value = 1413.5927
$ARM.Range("M61").Value = -1201.5927
$ARM.Range("H63").Value = 1305.5
$ARM.Range("I63").Value = 1111
$ARM.Range("J63").Value = 1500
$ARM.Range("K63").Value = 1111
$ARM.Range("L63").Value = 1500
$ARM.Range("M63").Value = -425
$ARM.Range("N63").Value = -2872
$ARM.Range("H66").Value = 1305.5
$ARM.Range("I66").Value = 1111
$ARM.Range("J66").Value = 1500
$ARM.Range("K66").Value = 5555
$ARM.Range("L66").Value = 7500
$ARM.Range("M66").Value = -2123
$ARM.Range("N66").Value = -14364
$ARM.Range("H88").Value = 8117
$ARM.Range("I88").Value = 10606
$ARM.Range("J88").Value = 7121.4
$ARM.Range("K88").Value = 10606
$ARM.Range("L88").Value = 7121.4
$ARM.Range("M88").Value = -10200
$ARM.Range("N88").Value = -7933.4
$ARM.Range("H91").Value = 8117
$ARM.Range("I91").Value = 10606
$ARM.Range("J91").Value = 7121.4
$ARM.Range("K91").Value = 10606
$ARM.Range("L91").Value = 7121.4
$ARM.Range("M91").Value = -9202
$ARM.Range("N91").Value = -9929.4
$ARM.Range("H99").Value = 11721.5
$ARM.Range("I99").Value = 754.4
$ARM.Range("J99").Value = 30000
$ARM.Range("K99").Value = 754.4
$ARM.Range("L99").Value = 30000
$ARM.Range("M99").Value = 2240.6
$ARM.Range("N99").Value = -35990
$ARM.Range("H132").Value = 6826.595
$ARM.Range("I132").Value = 7511.4707
$ARM.Range("J132").Value = 3915.875
$ARM.Range("K132").Value = 22534.4121
$ARM.Range("L132").Value = 11747.625
$ARM.Range("M132").Value = -20004.4121
$ARM.Range("N132").Value = -16807.625
$ARM.Range("H136").Value = 1702.5526
$ARM.Range("I136").Value = 1413.5927
$ARM.Range("K136").Value = 4240.7781
$ARM.Range("M136").Value = -1690.7781

# ---- BSM ----
$BSM = $wb.Worksheets.Item("BSM")
$BSM.Range("H86").Value = 2666.6667
$BSM.Range("I86").Value = 2562.5
$BSM.Range("J86").Value = 3500
$BSM.Range("K86").Value = 2562.5
$BSM.Range("L86").Value = 3500
$BSM.Range("M86").Value = -1439.5
$BSM.Range("N86").Value = -5746
$BSM.Range("H89").Value = 2666.6667
$BSM.Range("I89").Value = 2562.5
$BSM.Range("J89").Value = 3500
$BSM.Range("K89").Value = 12812.5
$BSM.Range("L89").Value = 17500
$BSM.Range("M89").Value = -7196.5
$BSM.Range("N89").Value = -28732
$BSM.Range("H105").Value = 2180
$BSM.Range("I105").Value = 2250
$BSM.Range("J105").Value = 1900
$BSM.Range("K105").Value = 2250
$BSM.Range("L105").Value = 1900
$BSM.Range("M105").Value = -503
$BSM.Range("N105").Value = -5394
$BSM.Range("H134").Value = 4396.628
$BSM.Range("I134").Value = 5031.129
$BSM.Range("K134").Value = 15093.387
$BSM.Range("M134").Value = -12558.387

# ---- CRP ----
$CRP = $wb.Worksheets.Item("CRP")
$CRP.Range("H16").Value = 804.5263
$CRP.Range("I16").Value = 759
$CRP.Range("J16").Value = 845.5
$CRP.Range("K16").Value = 759
$CRP.Range("L16").Value = 845.5
$CRP.Range("M16").Value = -472
$CRP.Range("N16").Value = -1419.5
$CRP.Range("H22").Value = 289.4
$CRP.Range("I22").Value = 257
$CRP.Range("J22").Value = 500
$CRP.Range("K22").Value = 257
$CRP.Range("L22").Value = 500
$CRP.Range("M22").Value = 93
$CRP.Range("N22").Value = -1200
$CRP.Range("H58").Value = 1312.9722
$CRP.Range("I58").Value = 1359.8788
$CRP.Range("J58").Value = 797
$CRP.Range("K58").Value = 1359.8788
$CRP.Range("L58").Value = 797
$CRP.Range("M58").Value = -1156.8788
$CRP.Range("N58").Value = -1203
$CRP.Range("H62").Value = 23812398
$CRP.Range("I62").Value = 2976.5454
$CRP.Range("K62").Value = 2976.5454
$CRP.Range("M62").Value = -2352.5454
$CRP.Range("H65").Value = 23812398
$CRP.Range("I65").Value = 2976.5454
$CRP.Range("K65").Value = 14882.727
$CRP.Range("M65").Value = -11762.727
$CRP.Range("H105").Value = 614.8
$CRP.Range("I105").Value = 614.8
$CRP.Range("J105").Value = 0
$CRP.Range("K105").Value = 614.8
$CRP.Range("L105").Value = 0
$CRP.Range("M105").Value = 1132.2
$CRP.Range("N105").ClearContents()
$CRP.Range("H113").Value = 804.5263
$CRP.Range("I113").Value = 759
$CRP.Range("J113").Value = 845.5
$CRP.Range("K113").Value = 759
$CRP.Range("L113").Value = 845.5
$CRP.Range("M113").Value = 1411
$CRP.Range("N113").Value = -5185.5
$CRP.Range("H134").Value = 5748.143
$CRP.Range("I134").Value = 6470.6113
$CRP.Range("J134").Value = 1413.3334
$CRP.Range("K134").Value = 19411.8339
$CRP.Range("L134").Value = 4240.0002
$CRP.Range("M134").Value = -16876.8339
$CRP.Range("N134").Value = -9310.0002
$CRP.Range("H136").Value = 1312.9722
$CRP.Range("I136").Value = 1359.8788
$CRP.Range("J136").Value = 797
$CRP.Range("K136").Value = 4079.6364
$CRP.Range("L136").Value = 2391
$CRP.Range("M136").Value = -1529.6364
$CRP.Range("N136").Value = -7491

# ---- CUL ----
$CUL = $wb.Worksheets.Item("CUL")
$CUL.Range("H2").Value = 257.88235
$CUL.Range("I2").Value = 319.25925
$CUL.Range("J2").Value = 21.142857
$CUL.Range("K2").Value = 1915.5555
$CUL.Range("L2").Value = 126.857142
$CUL.Range("M2").Value = -1802.5555
$CUL.Range("N2").Value = -352.857142
$CUL.Range("H140").Value = 2350.8286
$CUL.Range("I140").Value = 1399.3529
$CUL.Range("K140").Value = 4198.0587
$CUL.Range("M140").Value = 981.9412999999995

# ---- GSM ----
$GSM = $wb.Worksheets.Item("GSM")
$GSM.Range("H70").Value = 55559896
$GSM.Range("I70").Value = 78951620
$GSM.Range("J70").Value = 4564.75
$GSM.Range("K70").Value = 78951620
$GSM.Range("L70").Value = 4564.75
$GSM.Range("M70").Value = -78951350
$GSM.Range("N70").Value = -5104.75
$GSM.Range("H73").Value = 55559896
$GSM.Range("I73").Value = 78951620
$GSM.Range("J73").Value = 4564.75
$GSM.Range("K73").Value = 78951620
$GSM.Range("L73").Value = 4564.75
$GSM.Range("M73").Value = -78950684
$GSM.Range("N73").Value = -6436.75
$GSM.Range("H80").Value = 4554.591
$GSM.Range("I80").Value = 4629.5
$GSM.Range("K80").Value = 4629.5
$GSM.Range("M80").Value = -3631.5
$GSM.Range("H83").Value = 4554.591
$GSM.Range("I83").Value = 4629.5
$GSM.Range("K83").Value = 23147.5
$GSM.Range("M83").Value = -18155.5
$GSM.Range("H113").Value = 16667852
$GSM.Range("I113").Value = 50000920
$GSM.Range("J113").Value = 1317.9
$GSM.Range("K113").Value = 50000920
$GSM.Range("L113").Value = 1317.9
$GSM.Range("M113").Value = -49998750
$GSM.Range("N113").Value = -5657.9
$GSM.Range("H126").Value = 1929.1428
$GSM.Range("J126").Value = 2875
$GSM.Range("L126").Value = 8625
$GSM.Range("N126").Value = -13565
$GSM.Range("H132").Value = 5557.724
$GSM.Range("I132").Value = 5944.8335
$GSM.Range("K132").Value = 17834.5005
$GSM.Range("M132").Value = -15304.5005

# ---- LTW ----
$LTW = $wb.Worksheets.Item("LTW")
$LTW.Range("H122").Value = 11650.75
$LTW.Range("I122").Value = 11650.75
$LTW.Range("J122").Value = 0
$LTW.Range("K122").Value = 34952.25
$LTW.Range("L122").Value = 0
$LTW.Range("M122").Value = -32502.25
$LTW.Range("N122").ClearContents()

# ---- WVR ----
$WVR = $wb.Worksheets.Item("WVR")
$WVR.Range("H122").Value = 1727.963
$WVR.Range("I122").Value = 1610.4103
$WVR.Range("J122").Value = 2033.6
$WVR.Range("K122").Value = 4831.2309
$WVR.Range("L122").Value = 6100.799999999999
$WVR.Range("M122").Value = -2381.2309
$WVR.Range("N122").Value = -11000.8
